$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Chars" table in columns C/D (plus a styled literal-text cell in E2) ---
# Values are entered in this specific order so the shared-string table is
# appended in the same sequence as the target workbook.
$ws.Range("C2").Value = "é"
$ws.Range("C3").Value = "ñ"
$ws.Range("C4").Value = "Σ"
$ws.Range("C1").Value = "Chars"
$ws.Range("D1").Value = "Chars comment"
$ws.Range("D4").Value = "Greek Capital Letter Sigma cannot be encoded in latin-1"
$ws.Range("D2").Value = "Latin small letter e with acute"
$ws.Range("E2").Value = "&#233;"

# E2 gets a small Verdana font (built via a scratch named style so the font
# table only grows by the properties actually changed, then the style is
# removed again so cellStyles/cellStyleXfs stay at their original counts).
$charStyle = $wb.Styles.Add("CharsTempStyle")
$charStyle.Font.Name = "Verdana"
$charStyle.Font.Size = 9
$charStyle.Font.Color = 0
$ws.Range("E2").Style = "CharsTempStyle"
$wb.Styles("CharsTempStyle").Delete()

# Column D needs to be wide enough to show the longest comment.
$ws.Columns.Item(4).ColumnWidth = 50.451822916666664

# --- Totals row ---
$ws.Range("B6").Formula = "=SUM(B2:B4)"

# --- Page setup + selection, matching the edited workbook's view state ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
[void]$ws.Range("B7").Select()
